$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 17
$ws.Range("A17").Value = 181117
$ws.Range("D17").Value = "复习视频内容"
$ws.Range("H17").Value = "数据库使用"
$ws.Range("L17").Value = "不操作遗忘"

# Row 19
$ws.Range("A19").Value = 181118
$ws.Range("D19").Value = "数据库安装配置"
$ws.Range("H19").Value = "复习MySQL"
$ws.Range("L19").Value = "mybatis没操作过"

# Row 21
$ws.Range("A21").Value = 181119
$ws.Range("D21").Value = "看视频学习springboot"
$ws.Range("H21").Value = "准备学习mybatis"
$ws.Range("L21").Value = "配置出错"

# Row 23
$ws.Range("A23").Value = 181120
$ws.Range("D23").Value = "补充架构说明书"
$ws.Range("H23").Value = "整合Mybatis "
$ws.Range("L23").Value = "逻辑结构不清晰"

# Update the active selection to match the new focus area
$ws.Range("D27:G28").Select()
